$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    $ws.Range("E2").Value = 3
    $ws.Range("G2").Value = 51.59157666666666
    $ws.Range("H2").Value = 154.77473
    $ws.Range("I2").Value = 0.2641250550177587
    $ws.Range("J2").Value = 0.2641250550177588
    $ws.Range("K2").Value = 3
    $ws.Range("M2").Value = 34.071057
    $ws.Range("N2").Value = 102.213171
    $ws.Range("O2").Value = 0.5537562116045693
    $ws.Range("P2").Value = 0.5537562116045693
    $ws.Range("Q2").Value = 1757.77954932987
    $ws.Range("R2").Value = 15820.01594396883
    $ws.Range("S2").Value = 0.1462608898564825
    $ws.Range("T2").Value = 0.1462608898564826
    $ws.Range("E3").Value = 3
    $ws.Range("G3").Value = 51.59157666666666
    $ws.Range("H3").Value = 154.77473
    $ws.Range("I3").Value = 0.2641250550177587
    $ws.Range("J3").Value = 0.2641250550177588
    $ws.Range("K3").Value = 3
    $ws.Range("M3").Value = 19.28977566666667
    $ws.Range("N3").Value = 57.869327
    $ws.Range("O3").Value = 0.3135163401556734
    $ws.Range("P3").Value = 0.3135163401556735
    $ws.Range("Q3").Value = 995.1899401896344
    $ws.Range("R3").Value = 8956.709461706709
    $ws.Range("S3").Value = 0.0828075205925836
    $ws.Range("T3").Value = 0.08280752059258364
    $ws.Range("E4").Value = 3
    $ws.Range("G4").Value = 51.59157666666666
    $ws.Range("H4").Value = 154.77473
    $ws.Range("I4").Value = 0.2641250550177587
    $ws.Range("J4").Value = 0.2641250550177588
    $ws.Range("K4").Value = 3
    $ws.Range("M4").Value = 6.368545333333333
    $ws.Range("N4").Value = 19.105636
    $ws.Range("O4").Value = 0.1035078406055505
    $ws.Range("P4").Value = 0.1035078406055505
    $ws.Range("Q4").Value = 328.5632948198088
    $ws.Range("R4").Value = 2957.069653378279
    $ws.Range("S4").Value = 0.02733901409471043
    $ws.Range("T4").Value = 0.02733901409471043
    $ws.Range("E5").Value = 3
    $ws.Range("G5").Value = 51.59157666666666
    $ws.Range("H5").Value = 154.77473
    $ws.Range("I5").Value = 0.2641250550177587
    $ws.Range("J5").Value = 0.2641250550177588
    $ws.Range("K5").Value = 3
    $ws.Range("M5").Value = 1.7978
    $ws.Range("N5").Value = 5.3934
    $ws.Range("O5").Value = 0.02921960763420679
    $ws.Range("P5").Value = 0.02921960763420679
    $ws.Range("Q5").Value = 92.75133653133332
    $ws.Range("R5").Value = 834.7620287819998
    $ws.Range("S5").Value = 0.007717630473982191
    $ws.Range("T5").Value = 0.007717630473982194
    $ws.Range("E6").Value = 3
    $ws.Range("G6").Value = 19.32115333333334
    $ws.Range("H6").Value = 57.96346000000001
    $ws.Range("I6").Value = 0.09891538535728452
    $ws.Range("J6").Value = 0.09891538535728453
    $ws.Range("K6").Value = 3
    $ws.Range("M6").Value = 34.071057
    $ws.Range("N6").Value = 102.213171
    $ws.Range("O6").Value = 0.5537562116045693
    $ws.Range("P6").Value = 0.5537562116045693
    $ws.Range("Q6").Value = 658.2921165257402
    $ws.Range("R6").Value = 5924.629048731661
    $ws.Range("S6").Value = 0.05477500906485596
    $ws.Range("T6").Value = 0.05477500906485597
    $ws.Range("E7").Value = 3
    $ws.Range("G7").Value = 19.32115333333334
    $ws.Range("H7").Value = 57.96346000000001
    $ws.Range("I7").Value = 0.09891538535728452
    $ws.Range("J7").Value = 0.09891538535728453
    $ws.Range("K7").Value = 3
    $ws.Range("M7").Value = 19.28977566666667
    $ws.Range("N7").Value = 57.869327
    $ws.Range("O7").Value = 0.3135163401556734
    $ws.Range("P7").Value = 0.3135163401556735
    $ws.Range("Q7").Value = 372.700713421269
    $ws.Range("R7").Value = 3354.306420791421
    $ws.Range("S7").Value = 0.03101158960230393
    $ws.Range("T7").Value = 0.03101158960230394
    $ws.Range("E8").Value = 3
    $ws.Range("G8").Value = 19.32115333333334
    $ws.Range("H8").Value = 57.96346000000001
    $ws.Range("I8").Value = 0.09891538535728452
    $ws.Range("J8").Value = 0.09891538535728453
    $ws.Range("K8").Value = 3
    $ws.Range("M8").Value = 6.368545333333333
    $ws.Range("N8").Value = 19.105636
    $ws.Range("O8").Value = 0.1035078406055505
    $ws.Range("P8").Value = 0.1035078406055505
    $ws.Range("Q8").Value = 123.0476408956178
    $ws.Range("R8").Value = 1107.42876806056
    $ws.Range("S8").Value = 0.01023851794099841
    $ws.Range("T8").Value = 0.01023851794099841
    $ws.Range("E9").Value = 3
    $ws.Range("G9").Value = 19.32115333333334
    $ws.Range("H9").Value = 57.96346000000001
    $ws.Range("I9").Value = 0.09891538535728452
    $ws.Range("J9").Value = 0.09891538535728453
    $ws.Range("K9").Value = 3
    $ws.Range("M9").Value = 1.7978
    $ws.Range("N9").Value = 5.3934
    $ws.Range("O9").Value = 0.02921960763420679
    $ws.Range("P9").Value = 0.02921960763420679
    $ws.Range("Q9").Value = 34.73556946266667
    $ws.Range("R9").Value = 312.6201251640001
    $ws.Range("S9").Value = 0.002890268749126217
    $ws.Range("T9").Value = 0.002890268749126218
    $ws.Range("E10").Value = 3
    $ws.Range("G10").Value = 112.3724673333333
    $ws.Range("H10").Value = 337.117402
    $ws.Range("I10").Value = 0.5752951554216499
    $ws.Range("J10").Value = 0.57529515542165
    $ws.Range("K10").Value = 3
    $ws.Range("M10").Value = 34.071057
    $ws.Range("N10").Value = 102.213171
    $ws.Range("O10").Value = 0.5537562116045693
    $ws.Range("P10").Value = 0.5537562116045693
    $ws.Range("Q10").Value = 3828.648739744638
    $ws.Range("R10").Value = 34457.83865770174
    $ws.Range("S10").Value = 0.3185732658207547
    $ws.Range("T10").Value = 0.3185732658207548
    $ws.Range("E11").Value = 3
    $ws.Range("G11").Value = 112.3724673333333
    $ws.Range("H11").Value = 337.117402
    $ws.Range("I11").Value = 0.5752951554216499
    $ws.Range("J11").Value = 0.57529515542165
    $ws.Range("K11").Value = 3
    $ws.Range("M11").Value = 19.28977566666667
    $ws.Range("N11").Value = 57.869327
    $ws.Range("O11").Value = 0.3135163401556734
    $ws.Range("P11").Value = 0.3135163401556735
    $ws.Range("Q11").Value = 2167.639685969828
    $ws.Range("R11").Value = 19508.75717372845
    $ws.Range("S11").Value = 0.180364431637085
    $ws.Range("T11").Value = 0.1803644316370851
    $ws.Range("E12").Value = 3
    $ws.Range("G12").Value = 112.3724673333333
    $ws.Range("H12").Value = 337.117402
    $ws.Range("I12").Value = 0.5752951554216499
    $ws.Range("J12").Value = 0.57529515542165
    $ws.Range("K12").Value = 3
    $ws.Range("M12").Value = 6.368545333333333
    $ws.Range("N12").Value = 19.105636
    $ws.Range("O12").Value = 0.1035078406055505
    $ws.Range("P12").Value = 0.1035078406055505
    $ws.Range("Q12").Value = 715.6491524308523
    $ws.Range("R12").Value = 6440.84237187767
    $ws.Range("S12").Value = 0.05954755924852953
    $ws.Range("T12").Value = 0.05954755924852955
    $ws.Range("E13").Value = 3
    $ws.Range("G13").Value = 112.3724673333333
    $ws.Range("H13").Value = 337.117402
    $ws.Range("I13").Value = 0.5752951554216499
    $ws.Range("J13").Value = 0.57529515542165
    $ws.Range("K13").Value = 3
    $ws.Range("M13").Value = 1.7978
    $ws.Range("N13").Value = 5.3934
    $ws.Range("O13").Value = 0.02921960763420679
    $ws.Range("P13").Value = 0.02921960763420679
    $ws.Range("Q13").Value = 202.0232217718666
    $ws.Range("R13").Value = 1818.2089959468
    $ws.Range("S13").Value = 0.01680989871528062
    $ws.Range("T13").Value = 0.01680989871528063
    $ws.Range("E14").Value = 3
    $ws.Range("G14").Value = 12.044915
    $ws.Range("H14").Value = 36.134745
    $ws.Range("I14").Value = 0.06166440420330686
    $ws.Range("J14").Value = 0.06166440420330688
    $ws.Range("K14").Value = 3
    $ws.Range("M14").Value = 34.071057
    $ws.Range("N14").Value = 102.213171
    $ws.Range("O14").Value = 0.5537562116045693
    $ws.Range("P14").Value = 0.5537562116045693
    $ws.Range("Q14").Value = 410.3829855251551
    $ws.Range("R14").Value = 3693.446869726396
    $ws.Range("S14").Value = 0.03414704686247609
    $ws.Range("T14").Value = 0.0341470468624761
    $ws.Range("E15").Value = 3
    $ws.Range("G15").Value = 12.044915
    $ws.Range("H15").Value = 36.134745
    $ws.Range("I15").Value = 0.06166440420330686
    $ws.Range("J15").Value = 0.06166440420330688
    $ws.Range("K15").Value = 3
    $ws.Range("M15").Value = 19.28977566666667
    $ws.Range("N15").Value = 57.869327
    $ws.Range("O15").Value = 0.3135163401556734
    $ws.Range("P15").Value = 0.3135163401556735
    $ws.Range("Q15").Value = 232.3437082740684
    $ws.Range("R15").Value = 2091.093374466615
    $ws.Range("S15").Value = 0.01933279832370089
    $ws.Range("T15").Value = 0.0193327983237009
    $ws.Range("E16").Value = 3
    $ws.Range("G16").Value = 12.044915
    $ws.Range("H16").Value = 36.134745
    $ws.Range("I16").Value = 0.06166440420330686
    $ws.Range("J16").Value = 0.06166440420330688
    $ws.Range("K16").Value = 3
    $ws.Range("M16").Value = 6.368545333333333
    $ws.Range("N16").Value = 19.105636
    $ws.Range("O16").Value = 0.1035078406055505
    $ws.Range("P16").Value = 0.1035078406055505
    $ws.Range("Q16").Value = 76.70858721364667
    $ws.Range("R16").Value = 690.37728492282
    $ws.Range("S16").Value = 0.006382749321312125
    $ws.Range("T16").Value = 0.006382749321312126
    $ws.Range("E17").Value = 3
    $ws.Range("G17").Value = 12.044915
    $ws.Range("H17").Value = 36.134745
    $ws.Range("I17").Value = 0.06166440420330686
    $ws.Range("J17").Value = 0.06166440420330688
    $ws.Range("K17").Value = 3
    $ws.Range("M17").Value = 1.7978
    $ws.Range("N17").Value = 5.3934
    $ws.Range("O17").Value = 0.02921960763420679
    $ws.Range("P17").Value = 0.02921960763420679
    $ws.Range("Q17").Value = 21.654348187
    $ws.Range("R17").Value = 194.889133683
    $ws.Range("S17").Value = 0.001801809695817758
    $ws.Range("T17").Value = 0.001801809695817759
